# Apply the Tue Oct 24 07:20:51 UTC 2023 cryptos-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose refreshed "Price" text would otherwise be auto-parsed by
# Excel as a plain number (dropping the literal/trailing-zero text shown
# on the page) are pre-formatted as Text so the exact string sticks.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "33.917.01"
$ws.Range("E2").Value = "  +10.51%  "
$ws.Range("D3").Value = "1.809.55"
$ws.Range("E3").Value = "  +7.37%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "228.00"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "30.94"
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").Value = "47.30"
$ws.Range("E9").Value = "  +6.66%  "
$ws.Range("E10").Value = "  +5.46%  "
$ws.Range("D11").Value = "0.0664"
$ws.Range("E11").Value = "  +6.23%  "
$ws.Range("D12").Value = "0.0931"
$ws.Range("E12").Value = "  +2.78%  "
$ws.Range("D13").Value = "2.070.44"
$ws.Range("E13").Value = "  +7.47%  "
$ws.Range("D14").Value = "1.819.13"
$ws.Range("E14").Value = "  +7.79%  "
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "33.870.69"
$ws.Range("E16").Value = "  +10.30%  "
$ws.Range("D17").Value = "10.07"
$ws.Range("E17").Value = "  -3.74%  "
$ws.Range("E18").Value = "  +6.47%  "
$ws.Range("D19").Value = "69.14"
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("D20").Value = "255.09"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").Value = "0.0₃0740"
$ws.Range("E21").Value = "  +3.59%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "10.38"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").Value = "158.72"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "16.43"
$ws.Range("E27").Value = "  +3.74%  "
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("D29").Value = "7.02"
$ws.Range("E29").Value = "  +5.07%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "3.79"
$ws.Range("E31").Value = "  +9.01%  "
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("E33").Value = "  +5.55%  "
$ws.Range("E34").Value = "  +6.61%  "
$ws.Range("D35").Value = "1.542.20"
$ws.Range("E35").Value = "  +2.16%  "
$ws.Range("D36").Value = "1.79"
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("E38").Value = "  +4.29%  "
$ws.Range("D39").Value = "83.47"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  +5.45%  "
$ws.Range("E41").Value = "  +4.14%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").Value = "0.901"
$ws.Range("E43").Value = "  +7.14%  "
$ws.Range("E44").Value = "  +5.51%  "
$ws.Range("D45").Value = "0.0521"
$ws.Range("E45").Value = "  +4.27%  "
$ws.Range("E46").Value = "  +4.35%  "
$ws.Range("D47").Value = "1.955.06"
$ws.Range("E47").Value = "  +7.20%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("B49").Value = "MinaProtocolToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D49").Value = "1.02"
$ws.Range("E49").Value = "  +149.86%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "5.64"
$ws.Range("E50").Value = "  +4.13%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "52.22"
$ws.Range("E51").Value = "  +1.08%  "

Write-Host "Applied cryptos update"
